# Export PAD fitting result
#
# Adds 16 new gaussian-fit columns (X:AM) to the "data" sheet for the
# optimized beta parameters (and their errors) of the He/Ne PAD fits:
#   he_beta1_opt, he_beta2_opt, he_beta3_opt, he_beta4_opt,
#   he_beta1_err, he_beta2_err, he_beta3_err, he_beta4_err,
#   ne_beta1_opt, ne_beta2_opt, ne_beta3_opt, ne_beta4_opt,
#   ne_beta1_err, ne_beta2_err, ne_beta3_err, ne_beta4_err

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$headers = @(
    "he_beta1_opt", "he_beta2_opt", "he_beta3_opt", "he_beta4_opt",
    "he_beta1_err", "he_beta2_err", "he_beta3_err", "he_beta4_err",
    "ne_beta1_opt", "ne_beta2_opt", "ne_beta3_opt", "ne_beta4_opt",
    "ne_beta1_err", "ne_beta2_err", "ne_beta3_err", "ne_beta4_err"
)
$cols = @("X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM")

# Row 1 (header): new cells get the same bold / centered / thin-bordered
# style already used by the rest of the header row (B1:W1).
for ($i = 0; $i -lt $cols.Length; $i++) {
  $cell = $ws.Range($cols[$i] + "1")
  $cell.Value = $headers[$i]
  $cell.Font.Bold = $true
  $cell.HorizontalAlignment = -4108  # xlCenter
  $cell.VerticalAlignment = -4160    # xlTop
  $cell.Borders.LineStyle = 1
  $cell.Borders.Weight = 2
}

# Rows 2-12: fitted values, in row order (data row N -> phase bin N-2).
$rows = @(
    @(0.2303369152318839,0.8119733954076405,0.4632945935281905,0.9828243046141045,0.03521173857543195,0.04960701158519783,0.04895676373782606,0.05921308206521127,0.1896714777692846,-0.1304234423415812,0.04584238207094469,0.05911214877885158,0.01274594221447612,0.01668608901795437,0.01735101340441096,0.01983109752452807),
    @(-0.03012458819371898,0.8250563588640976,-0.06189235157483897,0.9582534448241412,0.03260334883820176,0.04658162870313017,0.04474281987802754,0.05526258960465341,0.05257815376905647,-0.1231086630957186,-0.08462712866495486,0.07257148802423746,0.01580123711355777,0.02083554278546272,0.02168421525778268,0.02477265084664302),
    @(-0.3100552968074747,0.8790588221112837,-0.5265608465781486,1.025681510265131,0.0260749832068067,0.03684694785947318,0.03616935251113901,0.0437280859434106,-0.1844533696583744,-0.1028888330137132,-0.05645210012185851,0.03003826795233839,0.01765159464892612,0.02309897062568117,0.02404265179175153,0.02746907436506632),
    @(-0.3817505671742484,0.8932486793363658,-0.3779513103387256,0.9300711809744323,0.03501410798132234,0.04908869605353904,0.04732125125637155,0.05725364864771489,-0.3157856298383602,-0.1138290711341478,0.08407562353972085,-0.001722233209699745,0.01688214172106668,0.021782236583616,0.02267569720129635,0.02589314619178982),
    @(-0.1789301034672548,0.8047356234518636,0.09228040506908995,0.976208744311826,0.04114317591171455,0.05815881889319935,0.05610071933556744,0.06944107755923816,-0.2491550709767125,-0.1199310031155369,0.2340087684605713,0.03829083742681549,0.01483039533407079,0.01929644147276807,0.02019885046802608,0.02293700298557831),
    @(0.1194011921865471,0.8705164201526875,0.5351846399503438,1.024155061315161,0.0333029135073506,0.04784159025676267,0.04708674135899893,0.05685462721387976,-0.02616051304666231,-0.1147170492852172,0.2929439845051983,0.04727254296278853,0.01696571212028748,0.02237573916783739,0.02351205143548536,0.02660336895172749),
    @(0.2898221083226984,0.7950592473770783,0.6709124312357876,1.039470983394997,0.04267509971700443,0.05953358931353882,0.06049759400005788,0.0719132286864952,0.1474232957356426,-0.1251682279656451,0.1833512926159607,0.07221634958947591,0.01418438477117493,0.01862511593359674,0.01944391160584168,0.02214296285325055),
    @(0.2095660992779006,0.9172839685127551,0.3763507309053348,1.088785708233849,0.04305288559505391,0.06199677647381171,0.0594772747513039,0.07377963931470817,0.179677026943932,-0.1385520602697954,0.001118194192661183,0.04545505040366168,0.01303327920338582,0.01708116411452447,0.01775228744296523,0.02029220588344172),
    @(-0.1324571673700495,0.8584230145577452,-0.2437985568431947,1.036059361707391,0.02870876315344539,0.04111737002686834,0.03951198870046486,0.04906789163759256,-0.04476612679253632,-0.150439566327175,-0.1064674594192832,0.03202639736560783,0.01362878949803541,0.01799149812440456,0.01871552431876865,0.02136209762344828),
    @(-0.3704204490349135,0.7421500673107432,-0.4556600951544652,1.036943461264211,0.03024086422810601,0.0413244627582933,0.04125154092293516,0.05035845253129417,-0.2529496981869608,-0.1667546729909404,-0.004575388411390203,-0.005532042724045593,0.02149274917820489,0.02800212487672342,0.02906995337308669,0.03322282162869005),
    @(-0.3117952286362083,0.778650669527331,-0.1862005118335083,1.012167361389565,0.03746964477947432,0.0519707585946367,0.05052220604255942,0.06268878667148603,-0.3147686271899816,-0.134772354096588,0.1396426286924229,-0.007470494473802007,0.01786959762221647,0.02307511361346789,0.02404114830182371,0.02741154826690088)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
  $rowVals = $rows[$r]
  $excelRow = $r + 2
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + $excelRow).Value = $rowVals[$i]
  }
}
